# Grades_sheet.xlsx update — "best accuracy till now"
# Updates student IDs (col A), per-question scores (col B), running totals
# (col C), and flags the now-blank/needs-review totals with a red fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Column C cells that are cleared and highlighted red (quiz needs
#    re-grading / was mis-scored). ClearContents + solid red interior.
# ---------------------------------------------------------------------
$redFillCells = @("C4", "C13", "C14", "C18")
foreach ($cellRef in $redFillCells) {
    $cell = $ws.Range($cellRef)
    $cell.ClearContents()
    $cell.Interior.Color = 255
}

# ---------------------------------------------------------------------
# 2) Column B cells that are simply cleared (no answer submitted / score
#    removed).
# ---------------------------------------------------------------------
$clearCells = @("B2", "B7", "B8", "B9", "B11", "B12", "B14")
foreach ($cellRef in $clearCells) {
    $ws.Range($cellRef).ClearContents()
}

# ---------------------------------------------------------------------
# 3) Cells receiving a new value that looks numeric but must stay text
#    (IDs / scores are stored as text throughout this sheet), so force
#    the Text number format before writing so Excel doesn't coerce the
#    literal into a number.
# ---------------------------------------------------------------------
$textForced = @{
    "B4"  = "3"
    "B5"  = "4"
    "A6"  = "1180274"
    "B6"  = "5"
    "A7"  = "1180056"
    "A9"  = "1180606"
    "A10" = "1180456"
    "B10" = "9"
    "A11" = "2200022"
    "A12" = "11180552"
    "A13" = "11180207"
    "B15" = "2"
    "B16" = "3"
    "A17" = "11170343"
    "B17" = "4"
    "A18" = "1180172"
    "B18" = "5"
}
foreach ($cellRef in $textForced.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $textForced[$cellRef]
}

# ---------------------------------------------------------------------
# 4) Plain text value (non-numeric, no coercion risk).
# ---------------------------------------------------------------------
$ws.Range("B13").Value = "Y"

# ---------------------------------------------------------------------
# 5) Column C cells that simply get a new numeric total.
# ---------------------------------------------------------------------
$numericValues = @{
    "C2"  = 2
    "C6"  = 2
    "C7"  = 2
    "C8"  = 5
    "C9"  = 5
    "C10" = 5
    "C11" = 3
    "C12" = 2
    "C17" = 2
}
foreach ($cellRef in $numericValues.Keys) {
    $ws.Range($cellRef).Value = $numericValues[$cellRef]
}
